$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to be treated as text, matching the source
# data which stores prices as formatted strings (e.g. thousands-dot
# separated values like "62.872.95") rather than numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.872.95"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.058.97"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.19"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.63"
$ws.Range("E6").Value = "  +3.28%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.053.98"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.24"
$ws.Range("E11").Value = "  +2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000221"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.32"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.558.06"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.112"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.873.78"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.062.95"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.59"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.74"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.693"
$ws.Range("E22").Value = "  -1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.99"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.26"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.03"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  -4.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.97"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  +4.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.86"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "58.76"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.30"
$ws.Range("E34").Value = "  -5.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("E35").Value = "  +7.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "478.72"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.241.20"
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0395"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0789"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.10"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.249"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  +5.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.02"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.99"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.108"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0517"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("E51").Value = "  +0.24%  "
